$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.885.47'
$ws.Range("E2").Value = '  +0.78%  '

# Row 3
$ws.Range("D3").Value = '3.329.90'
$ws.Range("E3").Value = '  +1.06%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.49'
$ws.Range("E5").Value = '  +0.74%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.15'
$ws.Range("E6").Value = '  +0.88%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +0.82%  '

# Row 9
$ws.Range("D9").Value = '3.325.17'
$ws.Range("E9").Value = '  +0.97%  '

# Row 10
$ws.Range("E10").Value = '  +5.02%  '

# Row 11
$ws.Range("E11").Value = '  +0.92%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.53'
$ws.Range("E12").Value = '  +4.16%  '

# Row 13
$ws.Range("E13").Value = '  +1.49%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '703.24'
$ws.Range("E14").Value = '  +2.11%  '

# Row 15
$ws.Range("D15").Value = '3.871.96'
$ws.Range("E15").Value = '  +1.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.46'
$ws.Range("E16").Value = '  +1.08%  '

# Row 17
$ws.Range("D17").Value = '67.923.27'
$ws.Range("E17").Value = '  +0.76%  '

# Row 18
$ws.Range("E18").Value = '  -0.65%  '

# Row 19
$ws.Range("D19").Value = '3.328.37'
$ws.Range("E19").Value = '  +0.87%  '

# Row 20
$ws.Range("E20").Value = '  +0.77%  '

# Row 21
$ws.Range("E21").Value = '  +2.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.897'
$ws.Range("E22").Value = '  +0.63%  '

# Row 23
$ws.Range("E23").Value = '  +3.96%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.09'
$ws.Range("E24").Value = '  +0.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.69'
$ws.Range("E25").Value = '  -0.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.91'
$ws.Range("E26").Value = '  +0.67%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.60'
$ws.Range("E28").Value = '  +3.75%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.12'
$ws.Range("E29").Value = '  -2.47%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.60'
$ws.Range("E30").Value = '  +1.89%  '

# Row 31
$ws.Range("E31").Value = '  +5.75%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '568.32'
$ws.Range("E32").Value = '  -0.59%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.03'
$ws.Range("E33").Value = '  +1.41%  '

# Row 34
$ws.Range("E34").Value = '  +2.79%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.36'
$ws.Range("E35").Value = '  +3.58%  '

# Row 36
$ws.Range("E36").Value = '  -0.01%  '

# Row 37
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.43'
$ws.Range("E37").Value = '  +3.74%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.682.57'
$ws.Range("E38").Value = '  -5.30%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.71'
$ws.Range("E39").Value = '  +9.27%  '

# Row 40
$ws.Range("E40").Value = '  +3.40%  '

# Row 41
$ws.Range("E41").Value = '  +2.57%  '

# Row 42
$ws.Range("E42").Value = '  +6.11%  '

# Row 43
$ws.Range("D43").Value = '0.0₃0678'
$ws.Range("E43").Value = '  +0.97%  '

# Row 44
$ws.Range("E44").Value = '  +2.60%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.34'
$ws.Range("E45").Value = '  -2.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0408'
$ws.Range("E46").Value = '  +0.46%  '

# Row 47
$ws.Range("E47").Value = '  +5.16%  '

# Row 48
$ws.Range("E48").Value = '  +0.92%  '

# Row 49
$ws.Range("E49").Value = '  -0.28%  '

# Row 50
$ws.Range("E50").Value = '  -2.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.19'
$ws.Range("E51").Value = '  +0.61%  '
